$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update file names (B column) and fix "c, s" -> "c,s" (D column) for rows 2-8
$ws.Range("B2").Value = "e1_321_1_1.jpeg"
$ws.Range("D2").Value = "c,s"

$ws.Range("B3").Value = "e2_321_1_2.jpeg"
$ws.Range("D3").Value = "c,s"

$ws.Range("B4").Value = "e3_321_1_3.jpeg"
$ws.Range("D4").Value = "c,s"

$ws.Range("B5").Value = "e4_321_2_2.jpeg"
$ws.Range("D5").Value = "c,s"

$ws.Range("B6").Value = "e5_321_2_1.jpeg"
$ws.Range("D6").Value = "c,s"

$ws.Range("B7").Value = "e6_321_2_0.jpeg"
$ws.Range("D7").Value = "c,s"

$ws.Range("B8").Value = "e7_321_3_1.jpeg"
$ws.Range("D8").Value = "c,s"

# Clear rows 9-29 entirely (all columns A-L), keeping rows' existing A/F formatting
$ws.Range("A9:L29").ClearContents()
